$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '93.436.92'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.96%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.463.05'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.94%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '625.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.43'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +10.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.392'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.46%  '

$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("E10").Value = '  +14.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.463.01'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +11.59%  '

$ws.Range("E13").Value = '  +6.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.14%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.116.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '93.319.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.89%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000249'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.74%  '

$ws.Range("E18").Value = '  +9.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.456.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +12.12%  '

$ws.Range("E22").Value = '  +19.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '503.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.96%  '

$ws.Range("E24").Value = '  +12.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.82'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.74%  '

$ws.Range("E26").Value = '  +5.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '91.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.15'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.641.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.57%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.140'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.82%  '

$ws.Range("B33").Value = 'Dai'
$ws.Range("C33").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.02%  '

$ws.Range("E34").Value = '  +1.00%  '

$ws.Range("E35").Value = '  +9.41%  '

$ws.Range("E36").Value = '  +11.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '29.50'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.87%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '572.06'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.41%  '

$ws.Range("E39").Value = '  +8.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.52%  '

$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.915'
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = '  +3.72%  '

$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '23.74'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.79%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0424'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +12.38%  '

$ws.Range("E46").Value = '  +5.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.53%  '

$ws.Range("E48").Value = '  +2.41%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.91%  '

$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.15'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.71%  '

$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.25%  '
